$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for the new columns.
#    Old layout:  ... L=apoio_medio | M=contribuicoes | N=media_contribuicoes | O=menor_ano | P=maior_ano
#    New layout:  ... L=apoio_medio | M=apoio_std | N=apoio_min | O=apoio_max | P=contribuicoes |
#                 Q=contribuicoes_med | R=contribuicoes_std | S=contribuicoes_min | T=contribuicoes_max |
#                 U=menor_ano | V=maior_ano
# ---------------------------------------------------------------------------
$ws.Range("M:O").EntireColumn.Insert()
$ws.Range("R:T").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2) Headers
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"

$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"

$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"

# ---------------------------------------------------------------------------
# 3) Data - apoio_medio (L) recalculated, plus the new apoio_* / contribuicoes_* columns
# ---------------------------------------------------------------------------
$apoioMedio = @(94.9905854649608, 110.6538302669828, 82.69413375987617, 91.3503645951285, 53.14416408875834)
$apoioStd   = @(47.6922463509549, 45.12744090621267, 30.90119243508478, 52.58131393014926, 7.563317519432532)
$apoioMin   = @(47.35034461927121, 39.22956647121969, 13.93896149503088, 21.61624650544615, 47.79609098250058)
$apoioMax   = @(305.2480444061168, 257.7853211115706, 194.2230576381307, 792.0360759681182, 58.4922371950161)

$contribStd = @(212.2582078460797, 378.1809353534696, 547.4955526904555, 401.4417134786221, 18.38477631085023)
$contribMin = @(35, 1, 3, 1, 10)
$contribMax = @(808, 1711, 5879, 6494, 36)

for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i

    $ws.Cells.Item($row, 12).Value = $apoioMedio[$i]   # L
    $ws.Cells.Item($row, 13).Value = $apoioStd[$i]     # M
    $ws.Cells.Item($row, 14).Value = $apoioMin[$i]     # N
    $ws.Cells.Item($row, 15).Value = $apoioMax[$i]     # O

    $ws.Cells.Item($row, 18).Value = $contribStd[$i]   # R
    $ws.Cells.Item($row, 19).Value = $contribMin[$i]   # S
    $ws.Cells.Item($row, 20).Value = $contribMax[$i]   # T
}

# ---------------------------------------------------------------------------
# 4) Number formats for the new columns (match the sibling columns' formats)
# ---------------------------------------------------------------------------
$ws.Range("M2:O6").NumberFormat = "R$ #,##0.00"
$ws.Range("R2:T6").NumberFormat = "#,##0"

Write-Host "edit applied"
